$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3623.8096
$ws.Range("I74").Value = 3768.75
$ws.Range("J74").Value = 3160
$ws.Range("K74").Value = 3768.75
$ws.Range("L74").Value = 3160
$ws.Range("M74").Value = -2832.75
$ws.Range("N74").Value = -5032
$ws.Range("H77").Value = 3623.8096
$ws.Range("I77").Value = 3768.75
$ws.Range("J77").Value = 3160
$ws.Range("K77").Value = 18843.75
$ws.Range("L77").Value = 15800
$ws.Range("M77").Value = -14163.75
$ws.Range("N77").Value = -25160
$ws.Range("H107").Value = 5729.65
$ws.Range("I107").Value = 6458.5293
$ws.Range("K107").Value = 6458.5293
$ws.Range("M107").Value = -4538.5293
$ws.Range("H111").Value = 2958
$ws.Range("I111").Value = 2933.3333
$ws.Range("J111").Value = 3032
$ws.Range("K111").Value = 8799.999899999999
$ws.Range("L111").Value = 9096
$ws.Range("M111").Value = -5732.999899999999
$ws.Range("N111").Value = -15230
$ws.Range("H115").Value = 4006
$ws.Range("I115").Value = 3815
$ws.Range("J115").Value = 4101.5
$ws.Range("K115").Value = 11445
$ws.Range("L115").Value = 12304.5
$ws.Range("M115").Value = -9878
$ws.Range("N115").Value = -15438.5
$ws.Range("H116").Value = 2538.0952
$ws.Range("I116").Value = 1692.3077
$ws.Range("J116").Value = 3912.5
$ws.Range("K116").Value = 1692.3077
$ws.Range("L116").Value = 3912.5
$ws.Range("M116").Value = 1749.6923
$ws.Range("N116").Value = -10796.5
$ws.Range("H134").Value = 41741.668
$ws.Range("J134").Value = 41741.668
$ws.Range("L134").Value = 41741.668
$ws.Range("N134").Value = -51881.668
$ws.Range("H135").Value = 35716476
$ws.Range("I135").Value = 1721.9166
$ws.Range("J135").Value = 250005000
$ws.Range("K135").Value = 15497.2494
$ws.Range("L135").Value = 2250045000
$ws.Range("M135").Value = -12962.2494
$ws.Range("N135").Value = -2250050070
$ws.Range("H136").Value = 44796.152
$ws.Range("J136").Value = 44796.152
$ws.Range("L136").Value = 44796.152
$ws.Range("N136").Value = -54996.152
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2244.4375
$ws.Range("I2").Value = 2334.0667
$ws.Range("J2").Value = 900
$ws.Range("K2").Value = 2334.0667
$ws.Range("L2").Value = 900
$ws.Range("M2").Value = -2221.0667
$ws.Range("N2").Value = -1126
$ws.Range("H32").Value = 9694.062
$ws.Range("I32").Value = 8577.949000000001
$ws.Range("K32").Value = 8577.949000000001
$ws.Range("M32").Value = -8290.949000000001
$ws.Range("H45").Value = 2386.2222
$ws.Range("I45").Value = 2288.4348
$ws.Range("J45").Value = 2948.5
$ws.Range("K45").Value = 2288.4348
$ws.Range("L45").Value = 2948.5
$ws.Range("M45").Value = -1911.4348
$ws.Range("N45").Value = -3702.5
$ws.Range("H110").Value = 1995.65
$ws.Range("I110").Value = 1947.3684
$ws.Range("K110").Value = 1947.3684
$ws.Range("M110").Value = 97.63159999999993
$ws.Range("H116").Value = 2244.4375
$ws.Range("I116").Value = 2334.0667
$ws.Range("J116").Value = 900
$ws.Range("K116").Value = 2334.0667
$ws.Range("L116").Value = 900
$ws.Range("M116").Value = -40.06669999999986
$ws.Range("N116").Value = -5488
$ws.Range("H122").Value = 2047.0435
$ws.Range("I122").Value = 1929.5294
$ws.Range("J122").Value = 2380
$ws.Range("K122").Value = 5788.5882
$ws.Range("L122").Value = 7140
$ws.Range("M122").Value = -3338.5882
$ws.Range("N122").Value = -12040
$ws.Range("H132").Value = 16669759
$ws.Range("I132").Value = 38463524
$ws.Range("J132").Value = 3938.5881
$ws.Range("K132").Value = 115390572
$ws.Range("L132").Value = 11815.7643
$ws.Range("M132").Value = -115388042
$ws.Range("N132").Value = -16875.7643
$ws.Range("H134").Value = 47145.445
$ws.Range("J134").Value = 47145.445
$ws.Range("L134").Value = 47145.445
$ws.Range("N134").Value = -57285.445
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2244.4375
$ws.Range("I3").Value = 2334.0667
$ws.Range("J3").Value = 900
$ws.Range("K3").Value = 2334.0667
$ws.Range("L3").Value = 900
$ws.Range("M3").Value = -2220.0667
$ws.Range("N3").Value = -1128
$ws.Range("H134").Value = 3303.9092
$ws.Range("I134").Value = 2197.8
$ws.Range("K134").Value = 6593.400000000001
$ws.Range("M134").Value = -4058.400000000001
$ws.Range("H140").Value = 31839
$ws.Range("J140").Value = 31839
$ws.Range("L140").Value = 31839
$ws.Range("N140").Value = -42199
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5957077
$ws.Range("I31").Value = 1815.8334
$ws.Range("J31").Value = 12828532
$ws.Range("K31").Value = 1815.8334
$ws.Range("L31").Value = 12828532
$ws.Range("M31").Value = -1520.8334
$ws.Range("N31").Value = -12829122
$ws.Range("H34").Value = 5957077
$ws.Range("I34").Value = 1815.8334
$ws.Range("J34").Value = 12828532
$ws.Range("K34").Value = 1815.8334
$ws.Range("L34").Value = 12828532
$ws.Range("M34").Value = -1613.8334
$ws.Range("N34").Value = -12828936
$ws.Range("H107").Value = 785
$ws.Range("I107").Value = 688.7143
$ws.Range("J107").Value = 934.7778
$ws.Range("K107").Value = 688.7143
$ws.Range("L107").Value = 934.7778
$ws.Range("M107").Value = 1231.2857
$ws.Range("N107").Value = -4774.7778
$ws.Range("H132").Value = 45316.562
$ws.Range("I132").Value = 1260.2084
$ws.Range("J132").Value = 177485.62
$ws.Range("K132").Value = 3780.6252
$ws.Range("L132").Value = 532456.86
$ws.Range("M132").Value = -1250.6252
$ws.Range("N132").Value = -537516.86
$ws.Range("H135").Value = 51060
$ws.Range("J135").Value = 51060
$ws.Range("L135").Value = 51060
$ws.Range("N135").Value = -61200
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5875.4443
$ws.Range("J3").Value = 10199.75
$ws.Range("L3").Value = 30599.25
$ws.Range("N3").Value = -30823.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 37217
$ws.Range("J93").Value = 37217
$ws.Range("L93").Value = 37217
$ws.Range("N93").Value = -40961
$ws.Range("H94").Value = 23990
$ws.Range("J94").Value = 23990
$ws.Range("L94").Value = 23990
$ws.Range("N94").Value = -25342
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492
$ws.Range("H107").Value = 4962.75
$ws.Range("I107").Value = 650.5
$ws.Range("J107").Value = 9275
$ws.Range("K107").Value = 650.5
$ws.Range("L107").Value = 9275
$ws.Range("M107").Value = 1269.5
$ws.Range("N107").Value = -13115
$ws.Range("H122").Value = 1759.6666
$ws.Range("I122").Value = 1756
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 5268
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -2818
$ws.Range("N122").Value = -10300
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8312.571
$ws.Range("I46").Value = 1108.8572
$ws.Range("J46").Value = 11914.429
$ws.Range("K46").Value = 1108.8572
$ws.Range("L46").Value = 11914.429
$ws.Range("M46").Value = -920.8571999999999
$ws.Range("N46").Value = -12290.429
$ws.Range("H55").Value = 512.4595
$ws.Range("I55").Value = 467.91666
$ws.Range("J55").Value = 594.6923
$ws.Range("K55").Value = 467.91666
$ws.Range("L55").Value = 594.6923
$ws.Range("M55").Value = -294.91666
$ws.Range("N55").Value = -940.6923
$ws.Range("H82").Value = 83333336
$ws.Range("J82").Value = 83333336
$ws.Range("L82").Value = 83333336
$ws.Range("N82").Value = -83334058
$ws.Range("H85").Value = 83333336
$ws.Range("J85").Value = 83333336
$ws.Range("L85").Value = 83333336
$ws.Range("N85").Value = -83335832
$ws.Range("H122").Value = 2107.6667
$ws.Range("I122").Value = 2127
$ws.Range("J122").Value = 1895
$ws.Range("K122").Value = 6381
$ws.Range("L122").Value = 5685
$ws.Range("M122").Value = -3931
$ws.Range("N122").Value = -10585
$ws.Range("H134").Value = 69659
$ws.Range("J134").Value = 69659
$ws.Range("L134").Value = 69659
$ws.Range("N134").Value = -79799
$ws.Range("H136").Value = 2705.1538
$ws.Range("I136").Value = 2055.3333
$ws.Range("J136").Value = 4167.25
$ws.Range("K136").Value = 6165.999899999999
$ws.Range("L136").Value = 12501.75
$ws.Range("M136").Value = -3615.999899999999
$ws.Range("N136").Value = -17601.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2000
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 2000
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H122").Value = 1500
$ws.Range("J122").Value = 1500
$ws.Range("L122").Value = 4500
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 1088624.1
$ws.Range("I132").Value = 1554032.2
$ws.Range("J132").Value = 2671.75
$ws.Range("K132").Value = 4662096.6
$ws.Range("L132").Value = 8015.25
$ws.Range("M132").Value = -4659566.6
$ws.Range("N132").Value = -13075.25
